$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Update "Latest Handoff Datetime" for the 4c1b3af7 row (row 4) in zh-cn and de-de sheets
$zhcn.Cells.Item(4, 8).Value = "2016-09-07 08:27:32"
$dede.Cells.Item(4, 8).Value = "2016-09-07 08:27:44"

# Update the aggregated "Latest HO Xliff Generate Date" for the same row (row 4) on Overview
$overview.Cells.Item(4, 7).Value = "2016-09-07 08:27:44"
